# Fruta / hortaliza, semanal
# Insert two new weekly price records (Naranja - Valencia, Primera & Segunda)
# at the top of this market's date-ordered block (row 307), pushing the
# existing rows 307-376 down to 309-378.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 307 (rows 307:308), shifting everything below down by 2.
$ws.Range("A307:A308").EntireRow.Insert()

# New row 307: Naranja, Valencia, Primera
$ws.Range("A307").Value = 4
$ws.Range("B307").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C307").Value = "Los Lagos"
$ws.Range("D307").Value = "2022-03-22"
$ws.Range("E307").Value = 10
$ws.Range("F307").Value = "Fruta"
$ws.Range("G307").Value = 100102
$ws.Range("H307").Value = "Cítricos"
$ws.Range("I307").Value = 100102005
$ws.Range("J307").Value = "Naranja"
$ws.Range("K307").Value = "Valencia"
$ws.Range("L307").Value = "Primera"
$ws.Range("M307").Value = 700
$ws.Range("N307").Value = 17000
$ws.Range("O307").Value = 18000
$ws.Range("P307").Value = 17500
$ws.Range("Q307").Value = "$/caja 15 kilos empedrada"
$ws.Range("R307").Value = "Región de O'Higgins"
$ws.Range("S307").Value = 1167
$ws.Range("T307").Value = 15

# New row 308: Naranja, Valencia, Segunda
$ws.Range("A308").Value = 4
$ws.Range("B308").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C308").Value = "Los Lagos"
$ws.Range("D308").Value = "2022-03-22"
$ws.Range("E308").Value = 10
$ws.Range("F308").Value = "Fruta"
$ws.Range("G308").Value = 100102
$ws.Range("H308").Value = "Cítricos"
$ws.Range("I308").Value = 100102005
$ws.Range("J308").Value = "Naranja"
$ws.Range("K308").Value = "Valencia"
$ws.Range("L308").Value = "Segunda"
$ws.Range("M308").Value = 350
$ws.Range("N308").Value = 15000
$ws.Range("O308").Value = 15000
$ws.Range("P308").Value = 15000
$ws.Range("Q308").Value = "$/caja 15 kilos empedrada"
$ws.Range("R308").Value = "Región de O'Higgins"
$ws.Range("S308").Value = 1000
$ws.Range("T308").Value = 15
